$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "...Despite their similarities, it is a very good idea to have
# both." -> "...Despite their similarities, it is a essential to have
# both." ("essential" becomes its own run)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("very good idea", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "essential"
}

# ---------------------------------------------------------------------------
# Edit 2: "...manipulation. See Image_QC.docx for more details." ->
# "...manipulation. See 02_Image_QC.docx for more details." with
# "02_Image_QC.docx" now bold.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("QCing raw data is a crucial step before doing any image manipulation. See Image_QC.docx for more details.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $rng.Text = "QCing raw data is a crucial step before doing any image manipulation. See 02_Image_QC.docx for more details."
}

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("02_Image_QC.docx for more details.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $boldLen = [int]"02_Image_QC.docx".Length
    $rng.MoveEnd(1, $boldLen - $rng.End + $rng.Start)
    $rng.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# Edit 3: "...preprocessing was completed as expected. See Image_QC.docx
# for more details." -> insert a new bold "02_" run right before the
# existing bold "Image_QC.docx" run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("completed as expected. See Image_QC.docx", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $prefixLen = [int]"completed as expected. See ".Length
    $rng.MoveStart(1, $prefixLen)
    $insPoint = $rng.Duplicate
    $insPoint.Collapse(1)
    $insPoint.InsertAfter("02_")
    $insPoint.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# Edit 4: "...throughout the rest of this documentation collective" ->
# "...throughout the rest of this documentation collective, most of which
# are located in 02_Extras/Image_Analysis_Guides"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute(" collective", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(", most of which are located in 02_Extras/Image_Analysis_Guides")
}
